# Daily "days remaining" refresh.
# Column D = total days, F = start date (yyyyMMdd, numeric), E = days remaining.
# For each data row: recompute remaining = (F + D days) - today.
#   - if remaining stays positive, just write remaining back to E (F unchanged).
#   - if remaining would be <= 0, the booking is renewed: F is reset to today
#     and E is reset to the full D (a fresh countdown starting today).
# "today" advanced by one day since the workbook was last refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = Get-Date -Year 2025 -Month 10 -Day 30 -Hour 0 -Minute 0 -Second 0
$todayOA = $today.ToOADate()

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fText = [string]([int64]$fVal)
    if ($fText.Length -ne 8) {
        # malformed start date (e.g. stray digit) - leave row untouched
        continue
    }

    $year = [int]$fText.Substring(0, 4)
    $month = [int]$fText.Substring(4, 2)
    $day = [int]$fText.Substring(6, 2)

    $startDate = Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0
    $endDate = $startDate.AddDays([double]$dVal)
    $remaining = [int]($endDate.ToOADate() - $todayOA)

    if ($remaining -le 0) {
        $eCell.Value2 = [int]$dVal
        $fCell.Value2 = [int]$today.ToString("yyyyMMdd")
    }
    else {
        $eCell.Value2 = $remaining
    }
}
